$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.809.48"
$ws.Range("E2").Value = "  -1.19%  "

# Row 3
$ws.Range("D3").Value = "'1.898.51"
$ws.Range("E3").Value = "  -0.71%  "

# Row 4
$ws.Range("E4").Value = "  +0.46%  "

# Row 5
$ws.Range("D5").Value = "'0.7658"
$ws.Range("E5").Value = "  +3.15%  "

# Row 6
$ws.Range("D6").Value = "'240.37"
$ws.Range("E6").Value = "  -1.36%  "

# Row 7
$ws.Range("E7").Value = "  +0.37%  "

# Row 8
$ws.Range("D8").Value = "'0.3045"
$ws.Range("E8").Value = "  -2.54%  "

# Row 9
$ws.Range("D9").Value = "'25.27"
$ws.Range("E9").Value = "  -5.43%  "

# Row 10
$ws.Range("D10").Value = "'0.06831"
$ws.Range("E10").Value = "  -1.86%  "

# Row 11
$ws.Range("D11").Value = "'0.07977"
$ws.Range("E11").Value = "  +0.14%  "

# Row 12
$ws.Range("D12").Value = "'1.887.03"
$ws.Range("E12").Value = "  -0.76%  "

# Row 13
$ws.Range("D13").Value = "'0.7359"
$ws.Range("E13").Value = "  -5.74%  "

# Row 14
$ws.Range("D14").Value = "'5.168"
$ws.Range("E14").Value = "  -1.92%  "

# Row 15
$ws.Range("D15").Value = "'91.00"
$ws.Range("E15").Value = "  -1.42%  "

# Row 16
$ws.Range("D16").Value = "'29.817.00"
$ws.Range("E16").Value = "  -1.14%  "

# Row 17
$ws.Range("D17").Value = "'13.75"
$ws.Range("E17").Value = "  -3.62%  "

# Row 18
$ws.Range("D18").Value = "'5.887"
$ws.Range("E18").Value = "  +0.43%  "

# Row 19
$ws.Range("D19").Value = "'245.25"
$ws.Range("E19").Value = "  +1.24%  "

# Row 20
$ws.Range("E20").Value = "  -1.70%  "

# Row 21
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
$ws.Range("D22").Value = "'2.131.01"
$ws.Range("E22").Value = "  -0.76%  "

# Row 23
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.52%  "

# Row 24
$ws.Range("D24").Value = "'6.891"
$ws.Range("E24").Value = "  -1.40%  "

# Row 25
$ws.Range("D25").Value = "'166.92"
$ws.Range("E25").Value = "  -0.33%  "

# Row 26
$ws.Range("D26").Value = "'9.240"
$ws.Range("E26").Value = "  -1.71%  "

# Row 27
$ws.Range("E27").Value = "  -2.23%  "

# Row 28
$ws.Range("D28").Value = "'0.1282"
$ws.Range("E28").Value = "  -0.03%  "

# Row 29
$ws.Range("D29").Value = "'2.029"
$ws.Range("E29").Value = "  -1.94%  "

# Row 30
$ws.Range("D30").Value = "'1.400"
$ws.Range("E30").Value = "  +4.52%  "

# Row 31
$ws.Range("D31").Value = "'1.512"
$ws.Range("E31").Value = "  -1.77%  "

# Row 32
$ws.Range("D32").Value = "'4.264"
$ws.Range("E32").Value = "  -1.60%  "

# Row 33
$ws.Range("D33").Value = "'4.070"
$ws.Range("E33").Value = "  -0.74%  "

# Row 34
$ws.Range("D34").Value = "'0.05265"
$ws.Range("E34").Value = "  +2.43%  "

# Row 35
$ws.Range("D35").Value = "'1.244"
$ws.Range("E35").Value = "  -4.13%  "

# Row 36
$ws.Range("D36").Value = "'0.7254"
$ws.Range("E36").Value = "  -2.29%  "

# Row 37
$ws.Range("D37").Value = "'2.717"
$ws.Range("E37").Value = "  +0.25%  "

# Row 38
$ws.Range("D38").Value = "'0.01911"
$ws.Range("E38").Value = "  -1.68%  "

# Row 39
$ws.Range("E39").Value = "  -0.72%  "

# Row 40
$ws.Range("D40").Value = "'6.188"
$ws.Range("E40").Value = "  -2.26%  "

# Row 41
$ws.Range("D41").Value = "'0.4398"
$ws.Range("E41").Value = "  -2.05%  "

# Row 42
$ws.Range("D42").Value = "'71.90"
$ws.Range("E42").Value = "  -4.01%  "

# Row 43
$ws.Range("E43").Value = "  +0.25%  "

# Row 44
$ws.Range("D44").Value = "'0.8329"
$ws.Range("E44").Value = "  -0.48%  "

# Row 45
$ws.Range("D45").Value = "'1.878"
$ws.Range("E45").Value = "  -3.94%  "

# Row 46
$ws.Range("D46").Value = "'7.597"
$ws.Range("E46").Value = "  -3.07%  "

# Row 47
$ws.Range("D47").Value = "'99.88"
$ws.Range("E47").Value = "  -1.33%  "

# Row 48
$ws.Range("D48").Value = "'9.752"
$ws.Range("E48").Value = "  -1.71%  "

# Row 49
$ws.Range("D49").Value = "'2.036.81"
$ws.Range("E49").Value = "  -0.46%  "

# Row 50
$ws.Range("D50").Value = "'36.08"
$ws.Range("E50").Value = "  -2.93%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05918"
$ws.Range("E51").Value = "  -1.06%  "
